# Weekend-Weekday model validation over different months
# Update evap, Inflow and Scalar sheets with new validation data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "evap" sheet (B4:B34) - replaced with the values that used to
#    live on the "Inflow" sheet.
# ---------------------------------------------------------------
$wsEvap = $wb.Worksheets.Item("evap")
$evapValues = @(1658,1655,1652,1650,1648,1646,1643,1640,1637,1634,1632,1630,1627,1625,1621,1618,1615,1614,1612,1609,1607,1605,1603,1601,1599,1597,1595,1592,1590,1588,1586)
$row = 4
foreach ($v in $evapValues) {
    $wsEvap.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# ---------------------------------------------------------------
# 2) "Inflow" sheet (B4:B34) - new inflow data for the month.
# ---------------------------------------------------------------
$wsInflow = $wb.Worksheets.Item("Inflow")
$inflowValues = @(3830,2818,3918,4446,2018,6232,3466,3338,3391,4572,1374,2797,1959,6463,3023,2551,2535,5622,3503,2680,6297,6402,5466,7411,2189,2186,6602,5950,5347,4347,8417)
$row = 4
foreach ($v in $inflowValues) {
    $wsInflow.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# ---------------------------------------------------------------
# 3) "Scalar" sheet - recomputed summary values reflecting the
#    updated Inflow/evap data above.
# ---------------------------------------------------------------
$wsScalar = $wb.Worksheets.Item("Scalar")
$wsScalar.Range("B12").Value = 11389586.429760002
$wsScalar.Range("C22").Value = 12354243.800000001
$wsScalar.Range("C25").Value = 11389586.429760002
$wsScalar.Range("B32").Value = -50229
$wsScalar.Range("C32").Value = -50229
$wsScalar.Range("D32").Value = -50229
$wsScalar.Range("C33").Value = 11389586.429760002
$wsScalar.Range("C34").Value = 11389586.429760002
$wsScalar.Range("B36").Value = 12354243.800000001
$wsScalar.Range("C36").Value = 12354243.800000001
$wsScalar.Range("D36").Value = 12354243.800000001
